# Updates cryptos list values (Price + Volume(1h) columns) to match the
# scraped data refresh from "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds scraped text (sometimes thousand-dot formatted,
# sometimes plain decimals with significant trailing zeros), so force the
# Text number format before assigning -- otherwise Excel auto-converts
# plain-decimal-looking strings (e.g. "83.00", "1.030") into real numbers
# and silently drops the trailing zeros / thousands-style formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.627.51"
$ws.Range("E2").Value = "  +2.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.97"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.035"
$ws.Range("E4").Value = "  +3.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.40"
$ws.Range("E5").Value = "  +3.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.030"
$ws.Range("E6").Value = "  +2.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4379"
$ws.Range("E7").Value = "  +1.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3752"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07406"
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8766"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.49"
$ws.Range("E11").Value = "  +2.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.863.69"
$ws.Range("E12").Value = "  -5.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.517"
$ws.Range("E13").Value = "  +3.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.702"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07196"
$ws.Range("E15").Value = "  +4.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.00"
$ws.Range("E16").Value = "  +3.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.036"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009039"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.030"
$ws.Range("E19").Value = "  +2.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.46"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.630.16"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.068.48"
$ws.Range("E24").Value = "  -4.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.63"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.940"
$ws.Range("E26").Value = "  +3.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.76"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.300"
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.940"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.39"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09081"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.209"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7680"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.516"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.884"
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.032"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.153"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01982"
$ws.Range("E38").Value = "  +3.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05282"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5184"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.813"
$ws.Range("E41").Value = "  +5.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1674"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.734"
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.586"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108.98"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.61"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.720"
$ws.Range("E47").Value = "  +4.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4661"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06399"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.888"
$ws.Range("E50").Value = "  +4.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.55"
$ws.Range("E51").Value = "  +5.88%  "

# Rows 22 and 23 only changed their Volume(1h) figure -- Price is unchanged.
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  +0.75%  "
